# Apply the crypto price/volume update described by the commit
# "Updated cryptos list on Wed Mar  1 07:57:24 UTC 2023 with GitHub Actions"
#
# This script updates the Price (D) and Volume(1h) (E) columns for most
# rows, and also fixes the ordering of the Toncoin / WrappedBTC rows
# (24 and 25), which were swapped along with their data.
#
# Note: several "Price" values look like plain numbers (e.g. "1.0000",
# "0.9997"). Excel would normally coerce such text into a numeric value,
# silently dropping the trailing zeros. To preserve the original text
# formatting (as in the source spreadsheet), those values are assigned
# with a leading apostrophe ('value), which is the standard Excel idiom
# for forcing text interpretation on an otherwise numeric-looking value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.724.56'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '1.651.72'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('D4').Value = '''1.0000'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''0.9999'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '''304.02'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').Value = '''0.3822'
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('D8').Value = '''51.37'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = '''0.3607'
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').Value = '''1.248'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').Value = '''0.08234'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').Value = '''1.000'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').Value = '''22.62'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = '''6.543'
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').Value = '''7.406'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '1.642.64'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '''97.29'
$ws.Range('E18').Value = '  +3.60%  '
$ws.Range('D19').Value = '''0.06964'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '''6.780'
$ws.Range('E20').Value = '  +4.15%  '
$ws.Range('D21').Value = '''17.71'
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('D22').Value = '''0.9997'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '''12.61'
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('B24').Value = 'WrappedBTC'
$ws.Range('C24').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D24').Value = '23.726.16'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''2.556'
$ws.Range('E25').Value = '  +4.24%  '
$ws.Range('D26').Value = '''3.072'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').Value = '''21.31'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '''151.72'
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('D29').Value = '''5.254'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').Value = '''135.01'
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('D31').Value = '1.835.43'
$ws.Range('E31').Value = '  +2.00%  '
$ws.Range('D32').Value = '''6.859'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('D33').Value = '''1.093'
$ws.Range('E33').Value = '  +6.06%  '
$ws.Range('D34').Value = '''11.89'
$ws.Range('E34').Value = '  +11.04%  '
$ws.Range('D35').Value = '''2.109'
$ws.Range('E35').Value = '  -5.71%  '
$ws.Range('D36').Value = '''0.02830'
$ws.Range('E36').Value = '  +2.94%  '
$ws.Range('D37').Value = '''0.2521'
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('D38').Value = '''0.08840'
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').Value = '''6.093'
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('D40').Value = '''0.07046'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('E41').Value = '  +6.07%  '
$ws.Range('D42').Value = '''0.7067'
$ws.Range('E42').Value = '  +1.11%  '
$ws.Range('D43').Value = '''1.337'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').Value = '''16.06'
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('D45').Value = '''0.6533'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = '''2.340'
$ws.Range('E46').Value = '  +2.87%  '
$ws.Range('D47').Value = '''0.9999'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').Value = '''3.982'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').Value = '''0.07986'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').Value = '''128.16'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('D51').Value = '''1.190'
$ws.Range('E51').Value = '  -0.22%  '
